# Improves parameterization using hashtable
# Switch the "runmode" flag from "n" to "y" on each test sheet, and leave
# the UI pointed at the testSuite sheet (matching the last-active state
# captured when the workbook was saved).

$wb = $excel.ActiveWorkbook

$wsAddCustomer = $wb.Worksheets.Item("AddCustomerTest")
$wsOpenAccount = $wb.Worksheets.Item("OpenAccountTest")
$wsTestSuite   = $wb.Worksheets.Item("testSuite")

# runmode: n -> y
$wsAddCustomer.Range("E5").Value = "y"
$wsOpenAccount.Range("C2").Value = "y"
$wsTestSuite.Range("B4").Value = "y"

# Restore each sheet's last selected cell
$wsAddCustomer.Range("E5").Select()
$wsOpenAccount.Range("E16").Select()

# testSuite becomes the active/selected tab, with B4 kept selected
$wsTestSuite.Activate()
$wsTestSuite.Range("B4").Select()
